$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that sat after the
#    "code in New Hampshire" run (it is not part of the Bookmarks
#    collection Word shows by default, but it is addressable by name).
$d.Bookmarks("_GoBack").Delete()

# 2. Append the new "group reflection" paragraphs at the very end of
#    the document body (before the final section break). Always
#    re-grab a fresh, collapsed range at the current end of the
#    document before each insertion so the range tracks the edits.

# -- empty paragraph --
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# -- scatter-plot paragraph --
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("Yet another way to think about visualizing this is to use a scatter plot to show the correlation between donations and per capita income levels.  On the X axis could be the donation levels per town, the Y axis could be the per capita level for that town and the dots could be color coded based on party.  ")

# -- empty paragraph --
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# -- "final option" paragraph --
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("A final option, suggested by Devin Shackle, would be a map like this showing the donation levels, but I’m not sure how I would correlate the per capita income levels in each town.")

# -- empty paragraph --
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# -- hyperlink paragraph --
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("http://www-personal.umich.edu/~mejn/cartograms/hiv1024x512.png")
$linkRange = $d.Range($d.Content.End - 64, $d.Content.End)
$d.Hyperlinks.Add($linkRange, "http://www-personal.umich.edu/~mejn/cartograms/hiv1024x512.png")

# -- empty paragraph that will hold the restored "_GoBack" bookmark --
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$d.Bookmarks.Add("_GoBack", $r)
